$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Incidental: A63 held the phone number as text; normalize it to a plain
# number (matches the upstream re-export of this row).
$ws.Cells.Item(63, 1).Value = 51616176

$row = 64

# Phone number is stored as text (leading apostrophe forces text entry,
# then reset the style so no stray number-format/quote-prefix survives).
$ws.Cells.Item($row, 1).Value = "'51616191"
$ws.Cells.Item($row, 1).Style = "Normal"

# amount (blank / not provided)
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "Cash"
$ws.Cells.Item($row, 4).Value = "2025-08-20T08:04:15"
$ws.Cells.Item($row, 5).Value = 120

# discount_applied (blank / not provided)
$ws.Cells.Item($row, 6).Value = "'"
$ws.Cells.Item($row, 6).Style = "Normal"

$ws.Cells.Item($row, 7).Value = 120
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
